$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 16:01"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6291776
$ws.Range("C4").Value = 1039
$ws.Range("D4").Value = 3547926
$ws.Range("E4").Value = 2553836
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 190014

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 3872670
$ws.Range("C6").Value = 23702
$ws.Range("D6").Value = 2984467
$ws.Range("E6").Value = 820534
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 183
$ws.Range("H6").Value = 67669

# Row 24: Irak
$ws.Range("A24").Value = "Irak"
$ws.Range("B24").Value = 247039
$ws.Range("C24").Value = 4755
$ws.Range("D24").Value = 187757
$ws.Range("E24").Value = 52007
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 7275

# Row 58: Ghana
$ws.Range("A58").Value = "Ghana"
$ws.Range("B58").Value = 44713
$ws.Range("C58").Value = 55
$ws.Range("D58").Value = 43577
$ws.Range("E58").Value = 856
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 4
$ws.Range("H58").Value = 280

# Row 62: Suiza
$ws.Range("A62").Value = "Suiza"
$ws.Range("B62").Value = 43127
$ws.Range("C62").Value = 364
$ws.Range("D62").Value = 36500
$ws.Range("E62").Value = 4614
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 2013

# Row 64: Uzbekistan
$ws.Range("A64").Value = "Uzbekistan"
$ws.Range("B64").Value = 42688
$ws.Range("C64").Value = 251
$ws.Range("D64").Value = 40081
$ws.Range("E64").Value = 2276
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 331

# Row 69: Serbia
$ws.Range("A69").Value = "Serbia"
$ws.Range("B69").Value = 31676
$ws.Range("C69").Value = 95
$ws.Range("D69").Value = 30188
$ws.Range("E69").Value = 770
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 718

# Row 76: Bosnia y Herzegovina
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 20804
$ws.Range("C76").Value = 287
$ws.Range("D76").Value = 14120
$ws.Range("E76").Value = 6048
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 9
$ws.Range("H76").Value = 636

# Row 77: Corea del Sur
$ws.Range("A77").Value = "Corea del Sur"
$ws.Range("B77").Value = 20644
$ws.Range("C77").Value = 195
$ws.Range("D77").Value = 15529
$ws.Range("E77").Value = 4786
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 329

# Row 86: Republica de Macedonia
$ws.Range("A86").Value = "Republica de Macedonia"
$ws.Range("B86").Value = 14762
$ws.Range("C86").Value = 162
$ws.Range("D86").Value = 11956
$ws.Range("E86").Value = 2200
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 606

# Row 89: Zambia
$ws.Range("A89").Value = "Zambia"
$ws.Range("B89").Value = 12523
$ws.Range("C89").Value = 108
$ws.Range("D89").Value = 11562
$ws.Range("E89").Value = 669
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 292

# Row 98: Tayikistan
$ws.Range("A98").Value = "Tayikistan"
$ws.Range("B98").Value = 8690
$ws.Range("C98").Value = 36
$ws.Range("D98").Value = 7482
$ws.Range("E98").Value = 1139
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 69

# Row 100: Haiti
$ws.Range("A100").Value = "Haiti"
$ws.Range("B100").Value = 8301
$ws.Range("C100").Value = 43
$ws.Range("D100").Value = 5870
$ws.Range("E100").Value = 2221
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 4
$ws.Range("H100").Value = 210

# Row 127: Uganda
$ws.Range("A127").Value = "Uganda"
$ws.Range("B127").Value = 3288
$ws.Range("C127").Value = 176
$ws.Range("D127").Value = 1532
$ws.Range("E127").Value = 1723
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 33

# Row 162: Birmania
$ws.Range("A162").Value = "Birmania"
$ws.Range("B162").Value = 1111
$ws.Range("C162").Value = 116
$ws.Range("D162").Value = 359
$ws.Range("E162").Value = 746
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 6

# Row 163: Belice
$ws.Range("A163").Value = "Belice"
$ws.Range("B163").Value = 1101
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 240
$ws.Range("E163").Value = 848
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 13

# Row 164: Lesoto
$ws.Range("A164").Value = "Lesoto"
$ws.Range("B164").Value = 1085
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 528
$ws.Range("E164").Value = 526
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 31

# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# Row 215: Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
